$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.5333333333333333
$ws.Range("C2").Value = 0.6153846153846154
$ws.Range("D2").Value = 0.5714285714285715
$ws.Range("E2").Value = 13
$ws.Range("B3").Value = 0.4444444444444444
$ws.Range("C3").Value = 0.3636363636363636
$ws.Range("D3").Value = 0.4
$ws.Range("E3").Value = 11
$ws.Range("B4").Value = 0.5
$ws.Range("C4").Value = 0.5
$ws.Range("D4").Value = 0.5
$ws.Range("E4").Value = 0.5
$ws.Range("B5").Value = 0.4888888888888889
$ws.Range("C5").Value = 0.4895104895104895
$ws.Range("D5").Value = 0.4857142857142858
$ws.Range("B6").Value = 0.4925925925925926
$ws.Range("C6").Value = 0.5
$ws.Range("D6").Value = 0.4928571428571429
$ws.Range("B7").Value = 0.5333333333333333
$ws.Range("C7").Value = 0.6153846153846154
$ws.Range("D7").Value = 0.5714285714285715
$ws.Range("E7").Value = 13
$ws.Range("B8").Value = 0.4444444444444444
$ws.Range("C8").Value = 0.3636363636363636
$ws.Range("D8").Value = 0.4
$ws.Range("E8").Value = 11
$ws.Range("B9").Value = 0.5
$ws.Range("C9").Value = 0.5
$ws.Range("D9").Value = 0.5
$ws.Range("E9").Value = 0.5
$ws.Range("B10").Value = 0.4888888888888889
$ws.Range("C10").Value = 0.4895104895104895
$ws.Range("D10").Value = 0.4857142857142858
$ws.Range("B11").Value = 0.4925925925925926
$ws.Range("C11").Value = 0.5
$ws.Range("D11").Value = 0.4928571428571429
$ws.Range("B12").Value = 0.5454545454545454
$ws.Range("C12").Value = 0.4615384615384616
$ws.Range("D12").Value = 0.4999999999999999
$ws.Range("E12").Value = 13
$ws.Range("B13").Value = 0.4615384615384616
$ws.Range("C13").Value = 0.5454545454545454
$ws.Range("D13").Value = 0.4999999999999999
$ws.Range("E13").Value = 11
$ws.Range("B14").Value = 0.5
$ws.Range("C14").Value = 0.5
$ws.Range("D14").Value = 0.5
$ws.Range("E14").Value = 0.5
$ws.Range("B15").Value = 0.5034965034965035
$ws.Range("C15").Value = 0.5034965034965035
$ws.Range("D15").Value = 0.4999999999999999
$ws.Range("B16").Value = 0.506993006993007
$ws.Range("C16").Value = 0.5
$ws.Range("D16").Value = 0.4999999999999998
$ws.Range("B17").Value = 0.5625
$ws.Range("C17").Value = 0.6923076923076923
$ws.Range("D17").Value = 0.6206896551724138
$ws.Range("E17").Value = 13
$ws.Range("B18").Value = 0.5
$ws.Range("C18").Value = 0.3636363636363636
$ws.Range("D18").Value = 0.4210526315789474
$ws.Range("E18").Value = 11
$ws.Range("B19").Value = 0.5416666666666666
$ws.Range("C19").Value = 0.5416666666666666
$ws.Range("D19").Value = 0.5416666666666666
$ws.Range("E19").Value = 0.5416666666666666
$ws.Range("B20").Value = 0.53125
$ws.Range("C20").Value = 0.5279720279720279
$ws.Range("D20").Value = 0.5208711433756806
$ws.Range("B21").Value = 0.5338541666666666
$ws.Range("C21").Value = 0.5416666666666666
$ws.Range("D21").Value = 0.529189352692075
$ws.Range("B22").Value = 0.5833333333333334
$ws.Range("C22").Value = 0.5384615384615384
$ws.Range("D22").Value = 0.5599999999999999
$ws.Range("E22").Value = 13
$ws.Range("B23").Value = 0.5
$ws.Range("C23").Value = 0.5454545454545454
$ws.Range("D23").Value = 0.5217391304347826
$ws.Range("E23").Value = 11
$ws.Range("B24").Value = 0.5416666666666666
$ws.Range("C24").Value = 0.5416666666666666
$ws.Range("D24").Value = 0.5416666666666666
$ws.Range("E24").Value = 0.5416666666666666
$ws.Range("B25").Value = 0.5416666666666667
$ws.Range("C25").Value = 0.5419580419580419
$ws.Range("D25").Value = 0.5408695652173913
$ws.Range("B26").Value = 0.545138888888889
$ws.Range("C26").Value = 0.5416666666666666
$ws.Range("D26").Value = 0.542463768115942
